# Time Sheet update: add new timesheet entries (rows 9-26), extend the
# D-column elapsed-time formula, widen column F for the longer notes,
# and move the selection/scroll to the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: finish time changed, a description was added -------------
$ws.Range("C9").Value = 0.82291666666666663
$ws.Range("E9").Value = "Program can now generate Team standings"

# --- New rows 10-26: Date / Start / End / Description / Note ----------
$rowData = @(
    @{ Row = 10; A = 41829;  B = 0.74305555555555547; C = 0.82638888888888884; E = "Learning Kivy" },
    @{ Row = 11; A = $null;  B = 0.85416666666666663; C = 0.89583333333333337; E = "Worked on file loading" },
    @{ Row = 12; A = 41832;  B = 0.43055555555555558; C = 0.51388888888888895; E = "Working on a useful Graphic User Interface" },
    @{ Row = 13; A = $null;  B = 0.65277777777777779; C = 0.73611111111111116; E = "Continued work on the GUI" },
    @{ Row = 14; A = 41833;  B = 0.375;               C = 0.4375;              E = "The gui can now load files and save them into a pdf"; F = "Began work on getting user input" },
    @{ Row = 15; A = $null;  B = 0.5625;              C = 0.60416666666666663; E = "The user can now enter division information" },
    @{ Row = 16; A = 41839;  B = 0.51041666666666663; C = 0.53125;             E = "Organized Project folders and files" },
    @{ Row = 17; A = 41846;  B = 0.47916666666666669; C = 0.5625;              E = "Get user input for tournament information" },
    @{ Row = 18; A = $null;  B = 0.61458333333333337; C = 0.64583333333333337; E = "Get highlight color choice" },
    @{ Row = 19; A = 41872;  B = 0.45833333333333331; C = 0.52083333333333337; E = "You can now select trophy winners highlight color" },
    @{ Row = 20; A = $null;  B = 0.5625;              C = 0.58333333333333337; E = "Setting up the screen for misc player identification" },
    @{ Row = 21; A = 41873;  B = 0.54166666666666663; C = 0.58333333333333337; E = "Adding functionality for misc player identification" },
    @{ Row = 22; A = 41876;  B = 0.4375;              C = 0.52083333333333337; E = "Allowed the selection of players for misc identification" },
    @{ Row = 23; A = 41880;  B = 0.41666666666666669; C = 0.45833333333333331; E = "Remove player identification" },
    @{ Row = 24; A = 41884;  B = 0.54166666666666663; C = 0.66666666666666663; E = "Saving the PDF"; F = "Severly Broke the project" },
    @{ Row = 25; A = $null;  B = 0.875;               C = 0.91666666666666663; E = "Fixing the project"; F = "Failed to fix the project" },
    @{ Row = 26; A = 41885;  B = 0.4375;              C = 0.47916666666666669; E = "Attempting to fix the way data is passed through out the application" }
)

# Apply the number formats (date / time) used elsewhere in the sheet to
# the new cells by copying the formatting from the existing row 2, then
# fill in the values. Only rows that actually have a date get a copy of
# the date format on column A, so blank-date rows keep no A cell at all
# (matching the rest of the sheet).
$ws.Range("A2").Copy()
foreach ($item in $rowData) {
    if ($item.A -ne $null) {
        $ws.Cells.Item($item.Row, 1).PasteSpecial(-4122)
    }
}

$ws.Range("B2:C2").Copy()
$ws.Range("B10:C26").PasteSpecial(-4122)

foreach ($item in $rowData) {
    $r = $item.Row
    if ($item.A -ne $null) { $ws.Cells.Item($r, 1).Value = $item.A }
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 5).Value = $item.E
    if ($item.ContainsKey("F")) {
        $ws.Cells.Item($r, 6).Value = $item.F
    }
}

# --- Column D: hours-elapsed formula for the new rows ------------------
$ws.Range("D22:D26").Formula = "=HOUR(C22-B22) + MINUTE(C22-B22) / 60"

# --- Column F: widen it to fit the longer notes -------------------------
$ws.Columns(6).ColumnWidth = 40.8

# --- Move the selection to the new bottom of the data (row 27) ---------
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A27").Select()
